# Applies:
#  1. Date placeholder text 4/26/2017 -> 10/23/2017 across the slide master,
#     all 11 slide layouts, and the notes master.
#  2. Slide 1 subtitle: autofit shrink + split last line / add "University of
#     Utah" line.
#  3. Slides 20 & 21: merge the "Level 1 (Participant Level):" and
#     "Level 2 (Group Level):" textbox runs back into a single run each.

$p = $ppt.ActivePresentation

function Set-DateFieldText($shape, [string]$newText) {
    $tr = $shape.TextFrame.TextRange
    # Assigning straight through leaves the text unchanged when the new
    # value happens to share a prefix with the old one (observed engine
    # quirk), so force a real change first.
    $tr.Text = "~"
    $tr.Text = $newText
}

# --- 1. Date placeholders -------------------------------------------------

$newDate = "10/23/2017"

# Slide master
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.Name -like "Date Placeholder*") {
        Set-DateFieldText $sh $newDate
    }
}

# Slide layouts
$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $sh = $layout.Shapes.Item($j)
        if ($sh.Name -like "Date Placeholder*") {
            Set-DateFieldText $sh $newDate
        }
    }
}

# Notes master
if ($p.HasNotesMaster) {
    $nm = $p.NotesMaster
} else {
    $nm = $p.NotesMaster
}
for ($i = 1; $i -le $nm.Shapes.Count; $i++) {
    $sh = $nm.Shapes.Item($i)
    if ($sh.Name -like "Date Placeholder*") {
        Set-DateFieldText $sh $newDate
    }
}

# --- 2. Slide 1 subtitle ----------------------------------------------------

$slide1 = $p.Slides.Item(1)
for ($i = 1; $i -le $slide1.Shapes.Count; $i++) {
    $sh = $slide1.Shapes.Item($i)
    if ($sh.Name -eq "Subtitle 2") {
        $subtitle = $sh
    }
}

$tr = $subtitle.TextFrame.TextRange
$full = $tr.Text
$idx = $full.IndexOf("Training")
if ($idx -ge 0) {
    # Force the trailing word into its own run (mirrors the author placing
    # the cursor there before continuing to type).
    $sub = $tr.Characters($idx + 1, 8)
    $sub.Text = "Training"
}

# Add the new "University of Utah" paragraph after the existing text.
$tr.InsertAfter("`rUniversity of Utah")

# --- 3. Slides 20 & 21: collapse multi-run "Level n (...)" labels ---------

function Merge-LevelLabel($shape, [string]$finalText) {
    $tr = $shape.TextFrame.TextRange
    # The current text is already split across several <a:r> runs; simply
    # re-assigning the same value is a no-op in this engine, so bounce
    # through a placeholder value first to force a genuine rewrite into a
    # single run.
    $tr.Text = "~placeholder~"
    $tr.Text = $finalText
}

foreach ($slideIdx in 20, 21) {
    $s = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $sh = $s.Shapes.Item($i)
        if ($sh.HasTextFrame -eq -1) {
            $t = $sh.TextFrame.TextRange.Text
            if ($t -eq "Level 1 (Participant Level):") {
                Merge-LevelLabel $sh "Level 1 (Participant Level):"
            } elseif ($t -eq "Level 2 (Group Level):") {
                Merge-LevelLabel $sh "Level 2 (Group Level):"
            }
        }
    }
}
